$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header labels (rearranged)
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# Data rows 2-6 (row 7 stays unchanged)
$row2 = @(0, 0, 0, 0, 1, 0)
$row3 = @(0, 0, 0, 0, 0, 1)
$row4 = @(1, 0, 0, 0, 0, 0)
$row5 = @(0, 1, 0, 0, 0, 0)
$row6 = @(0, 0, 1, 0, 0, 0)

for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item(2, $c).Value = $row2[$c - 1]
    $ws.Cells.Item(3, $c).Value = $row3[$c - 1]
    $ws.Cells.Item(4, $c).Value = $row4[$c - 1]
    $ws.Cells.Item(5, $c).Value = $row5[$c - 1]
    $ws.Cells.Item(6, $c).Value = $row6[$c - 1]
}
